$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# Row 1: header labels, Row 2: first vocabulary entry (word column only)
$ws.Cells.Item(1, 1).Value = "単語"
$ws.Cells.Item(1, 2).Value = "意味"
$ws.Cells.Item(2, 1).Value = "こうか"

# Touch row 3 so it becomes part of the used range, matching the
# A1:B3 layout of the uploaded word list (kept blank otherwise).
$ws.Range("A1:B3").Font.Name = "ＭＳ 明朝"
$ws.Range("A1:B3").Font.Size = 12

# --- Clean up the old bold/bordered header formatting -----------------
$ws.Range("A1:B1").Font.Bold = $false
$ws.Range("A1:B1").Borders.LineStyle = -4142

# --- Header alignment (centered horizontally, top vertically) ---------
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").VerticalAlignment = -4160

# --- Row heights --------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14
$ws.Rows.Item(2).RowHeight = 14
$ws.Rows.Item(3).RowHeight = 14

# --- Selection / active cell, matching the uploaded file ---------------
[void]$ws.Range("E8").Select()
